$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1987767584097859
$ws.Range("C2").Value = 0.5168195718654435
$ws.Range("J2").Value = 0.01834862385321101
$ws.Range("P2").Value = 0.1284403669724771
$ws.Range("S2").Value = 0.1376146788990826
$ws.Range("B3").Value = 0.01136363636363636
$ws.Range("C3").Value = 0.03977272727272727
$ws.Range("J3").Value = 0.01704545454545454
$ws.Range("P3").Value = 0.7272727272727273
$ws.Range("S3").Value = 0.2045454545454546
$ws.Range("J4").Value = 0.05660377358490566
$ws.Range("O4").Value = 0.01886792452830189
$ws.Range("P4").Value = 0.6415094339622641
$ws.Range("S4").Value = 0.2830188679245283
$ws.Range("B6").Value = 0.05179282868525897
$ws.Range("D6").Value = 0.00796812749003984
$ws.Range("F6").Value = 0.05577689243027888
$ws.Range("J6").Value = 0.2191235059760956
$ws.Range("O6").Value = 0.03187250996015936
$ws.Range("Q6").Value = 0.1673306772908366
$ws.Range("R6").Value = 0.06772908366533864
$ws.Range("S6").Value = 0.398406374501992
$ws.Range("B7").Value = 0.1272727272727273
$ws.Range("D7").Value = 0.01818181818181818
$ws.Range("F7").Value = 0.06363636363636363
$ws.Range("J7").Value = 0.1636363636363636
$ws.Range("O7").Value = 0.02272727272727273
$ws.Range("Q7").Value = 0.15
$ws.Range("R7").Value = 0.1
$ws.Range("S7").Value = 0.3545454545454546
$ws.Range("B8").Value = 0.1045751633986928
$ws.Range("D8").Value = 0.02396514161220044
$ws.Range("E8").Value = 0.002178649237472767
$ws.Range("F8").Value = 0.07843137254901961
$ws.Range("J8").Value = 0.1241830065359477
$ws.Range("O8").Value = 0.008714596949891068
$ws.Range("Q8").Value = 0.196078431372549
$ws.Range("R8").Value = 0.07625272331154684
$ws.Range("S8").Value = 0.3856209150326798
$ws.Range("B9").Value = 0.07834101382488479
$ws.Range("D9").Value = 0.01382488479262673
$ws.Range("E9").Value = 0.004608294930875576
$ws.Range("F9").Value = 0.07373271889400922
$ws.Range("J9").Value = 0.1336405529953917
$ws.Range("O9").Value = 0.04147465437788019
$ws.Range("Q9").Value = 0.1612903225806452
$ws.Range("R9").Value = 0.1059907834101382
$ws.Range("S9").Value = 0.3870967741935484
$ws.Range("B10").Value = 0.1037140854940434
$ws.Range("D10").Value = 0.02242466713384723
$ws.Range("F10").Value = 0.07778556412053259
$ws.Range("J10").Value = 0.1065171688857744
$ws.Range("O10").Value = 0.02522775052557814
$ws.Range("Q10").Value = 0.1983181499649615
$ws.Range("R10").Value = 0.08829712683952348
$ws.Range("S10").Value = 0.3777154870357393
$ws.Range("G11").Value = 0.1606648199445983
$ws.Range("J11").Value = 0.0997229916897507
$ws.Range("K11").Value = 0.2132963988919667
$ws.Range("L11").Value = 0.5069252077562327
$ws.Range("S11").Value = 0.01939058171745152
$ws.Range("G12").Value = 0.7377049180327869
$ws.Range("J12").Value = 0.2349726775956284
$ws.Range("K12").Value = 0.00546448087431694
$ws.Range("S12").Value = 0.02185792349726776
$ws.Range("G13").Value = 0.6415094339622641
$ws.Range("J13").Value = 0.2641509433962264
$ws.Range("S13").Value = 0.09433962264150944
$ws.Range("F15").Value = 0.01351351351351351
$ws.Range("H15").Value = 0.1576576576576577
$ws.Range("I15").Value = 0.03603603603603604
$ws.Range("J15").Value = 0.3378378378378378
$ws.Range("K15").Value = 0.07207207207207207
$ws.Range("M15").Value = 0.02252252252252252
$ws.Range("O15").Value = 0.04054054054054054
$ws.Range("S15").Value = 0.3198198198198198
$ws.Range("F16").Value = 0.02512562814070352
$ws.Range("H16").Value = 0.1457286432160804
$ws.Range("I16").Value = 0.07035175879396985
$ws.Range("J16").Value = 0.4522613065326633
$ws.Range("K16").Value = 0.1256281407035176
$ws.Range("M16").Value = 0.03015075376884422
$ws.Range("O16").Value = 0.01507537688442211
$ws.Range("S16").Value = 0.135678391959799
$ws.Range("F17").Value = 0.01464435146443515
$ws.Range("H17").Value = 0.1527196652719665
$ws.Range("I17").Value = 0.1129707112970711
$ws.Range("J17").Value = 0.4351464435146444
$ws.Range("K17").Value = 0.0899581589958159
$ws.Range("M17").Value = 0.02719665271966527
$ws.Range("N17").Value = 0.002092050209205021
$ws.Range("O17").Value = 0.04811715481171548
$ws.Range("S17").Value = 0.1171548117154812
$ws.Range("F18").Value = 0.01345291479820628
$ws.Range("H18").Value = 0.2152466367713005
$ws.Range("I18").Value = 0.1255605381165919
$ws.Range("J18").Value = 0.4035874439461883
$ws.Range("K18").Value = 0.08071748878923767
$ws.Range("M18").Value = 0.02242152466367713
$ws.Range("O18").Value = 0.03139013452914798
$ws.Range("S18").Value = 0.1076233183856502
$ws.Range("F19").Value = 0.01035911602209945
$ws.Range("H19").Value = 0.1926795580110497
$ws.Range("I19").Value = 0.08011049723756906
$ws.Range("J19").Value = 0.3791436464088398
$ws.Range("K19").Value = 0.1187845303867403
$ws.Range("M19").Value = 0.01933701657458563
$ws.Range("O19").Value = 0.05939226519337017
$ws.Range("S19").Value = 0.1401933701657458
